# Update slide_list path syntax:
#  - strip the leading ".\" from root/file path values (".\test_files" -> "test_files",
#    and the embedded tiff paths on the "pages" sheet)
#  - rename the "Sub" crop-region labels to "Crop" on the "slides" sheet
#  - re-activate the "slides" sheet/tab and update each sheet's selected cell

$wb = $excel.ActiveWorkbook

$wsSlides = $wb.Worksheets.Item("slides")
$wsPages  = $wb.Worksheets.Item("pages")

# --- "pages" sheet: filename column (G) path syntax ---
$wsPages.Range("G2").Value = "test_files\test_figure_output_1.tiff"
$wsPages.Range("G3").Value = "test_files\test_figure_output_2.tiff"

# --- "slides" sheet: root column (E) path syntax + Sub -> Crop labels (C) ---
$wsSlides.Range("E2").Value = "test_files"
$wsSlides.Range("E3").Value = "test_files"
$wsSlides.Range("C4").Value = "Slide1Crop"
$wsSlides.Range("E4").Value = "test_files"
$wsSlides.Range("C5").Value = "Slide2Crop"
$wsSlides.Range("E5").Value = "test_files"

# --- selection / active-tab bookkeeping ---
$wsPages.Range("D3").Select()

$wsSlides.Activate()
$wsSlides.Range("C5").Select()
